$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107 (shifts existing rows 107.. down by one,
# inheriting the formatting of the row above for the new row).
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new price-observation record.
$ws.Cells.Item(107, 1).Value  = 1
$ws.Cells.Item(107, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(107, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(107, 4).Value  = 44994
$ws.Cells.Item(107, 5).Value  = 15
$ws.Cells.Item(107, 6).Value  = "Fruta"
$ws.Cells.Item(107, 7).Value  = 100108
$ws.Cells.Item(107, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(107, 9).Value  = 100108002
$ws.Cells.Item(107, 10).Value = "Mango"
$ws.Cells.Item(107, 11).Value = "Sin especificar"
$ws.Cells.Item(107, 12).Value = "Especial"
$ws.Cells.Item(107, 13).Value = 456
$ws.Cells.Item(107, 14).Value = 4500
$ws.Cells.Item(107, 15).Value = 5000
$ws.Cells.Item(107, 16).Value = 4750
$ws.Cells.Item(107, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(107, 18).Value = "Perú"
$ws.Cells.Item(107, 19).Value = 1188
$ws.Cells.Item(107, 20).Value = 4
